$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.503.94"
$ws.Range("E2").Value = "'  -1.27%  "
$ws.Range("D3").Value = "'3.849.37"
$ws.Range("E3").Value = "'  -0.83%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'602.39"
$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("D6").Value = "'168.50"
$ws.Range("E6").Value = "'  -0.75%  "
$ws.Range("D7").Value = "'3.850.82"
$ws.Range("E7").Value = "'  -0.78%  "
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "'  -1.10%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "'  -0.73%  "
$ws.Range("D11").Value = "'6.50"
$ws.Range("E11").Value = "'  +1.84%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "'  -1.72%  "
$ws.Range("D13").Value = "'0.0000277"
$ws.Range("E13").Value = "'  +9.34%  "
$ws.Range("D14").Value = "'36.89"
$ws.Range("E14").Value = "'  -2.86%  "
$ws.Range("D15").Value = "'4.490.55"
$ws.Range("E15").Value = "'  -0.79%  "
$ws.Range("D16").Value = "'3.840.53"
$ws.Range("E16").Value = "'  -0.61%  "
$ws.Range("D17").Value = "'68.509.00"
$ws.Range("E17").Value = "'  -1.20%  "
$ws.Range("D18").Value = "'18.45"
$ws.Range("E18").Value = "'  -0.94%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "'  -3.35%  "
$ws.Range("E20").Value = "'  -1.05%  "
$ws.Range("D21").Value = "'10.93"
$ws.Range("E21").Value = "'  -0.16%  "
$ws.Range("D22").Value = "'470.96"
$ws.Range("E22").Value = "'  -3.23%  "
$ws.Range("D23").Value = "'0.731"
$ws.Range("E23").Value = "'  -1.85%  "
$ws.Range("D24").Value = "'0.0000161"
$ws.Range("E24").Value = "'  -4.48%  "
$ws.Range("D25").Value = "'83.33"
$ws.Range("E25").Value = "'  -2.04%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "'  -1.74%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "'  -1.92%  "
$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = "'  +1.85%  "
$ws.Range("E29").Value = "'  +0.03%  "
$ws.Range("E30").Value = "'  -0.83%  "
$ws.Range("D31").Value = "'3.995.83"
$ws.Range("E31").Value = "'  -0.74%  "
$ws.Range("D32").Value = "'7.66"
$ws.Range("E32").Value = "'  -2.69%  "
$ws.Range("D33").Value = "'31.31"
$ws.Range("E33").Value = "'  -1.11%  "
$ws.Range("E34").Value = "'  -2.01%  "
$ws.Range("D35").Value = "'9.31"
$ws.Range("E35").Value = "'  -2.76%  "
$ws.Range("D36").Value = "'3.812.12"
$ws.Range("E36").Value = "'  -0.80%  "
$ws.Range("B37").Value = "'Hedera"
$ws.Range("C37").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.104"
$ws.Range("E37").Value = "'  -1.45%  "
$ws.Range("B38").Value = "'dogwifhat"
$ws.Range("C38").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.73"
$ws.Range("E38").Value = "'  +11.17%  "
$ws.Range("B39").Value = "'Kaspa"
$ws.Range("C39").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.140"
$ws.Range("B40").Value = "'Mantle"
$ws.Range("C40").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.02"
$ws.Range("E40").Value = "'  -1.60%  "
$ws.Range("D41").Value = "'5.92"
$ws.Range("E41").Value = "'  -2.36%  "
$ws.Range("E42").Value = "'  +0.04%  "
$ws.Range("E43").Value = "'  -2.31%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("E44").Value = "'  -3.85%  "
$ws.Range("B45").Value = "'Bittensor"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'419.91"
$ws.Range("E45").Value = "'  -2.54%  "
$ws.Range("B46").Value = "'Cosmos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "'8.71"
$ws.Range("E46").Value = "'  +0.53%  "
$ws.Range("B47").Value = "'FLOKI"
$ws.Range("C47").Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000294"
$ws.Range("E47").Value = "'  +5.53%  "
$ws.Range("E48").Value = "'  -0.02%  "
$ws.Range("D49").Value = "'47.01"
$ws.Range("E49").Value = "'  -2.19%  "
$ws.Range("D50").Value = "'142.30"
$ws.Range("E50").Value = "'  +0.60%  "
$ws.Range("D51").Value = "'26.10"
$ws.Range("E51").Value = "'  +4.56%  "
